$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each Price/Volume cell is stored as text (inlineStr) in the workbook, so we
# force a Text number format before assigning the new value (otherwise Excel
# would reinterpret "298.99" or "1.60%" as numeric/percentage values), then
# reset the style back to Normal so no stray formatting is left behind.
function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "298.99"
Set-TextValue "E2" "1.60%"
Set-TextValue "D3" "32.08"
Set-TextValue "E3" "3.35%"
Set-TextValue "D4" "5.001"
Set-TextValue "E4" "1.45%"
Set-TextValue "D5" "0.07716"
Set-TextValue "E5" "5.17%"
Set-TextValue "D6" "2.262"
Set-TextValue "E6" "-1.45%"
Set-TextValue "D7" "7.916"
Set-TextValue "E7" "2.22%"
Set-TextValue "D8" "0.9220"
Set-TextValue "E8" "1.49%"
Set-TextValue "D9" "0.09927"
Set-TextValue "E9" "24.36%"
Set-TextValue "D10" "0.1753"
Set-TextValue "E10" "3.83%"
Set-TextValue "D11" "0.08411"
Set-TextValue "E11" "3.24%"
Set-TextValue "D12" "0.03323"
Set-TextValue "E12" "7.15%"
Set-TextValue "D13" "0.09869"
Set-TextValue "E13" "-2.15%"
Set-TextValue "D14" "0.001491"
Set-TextValue "E14" "-1.15%"
Set-TextValue "D15" "0.005701"
Set-TextValue "E15" "-2.71%"
Set-TextValue "D16" "3.536"
Set-TextValue "E16" "1.51%"
Set-TextValue "D17" "3.822"
Set-TextValue "E17" "2.04%"
Set-TextValue "E18" "4.90%"
Set-TextValue "D19" "0.3369"
Set-TextValue "E19" "1.22%"
Set-TextValue "D20" "0.1333"
Set-TextValue "E20" "2.14%"
Set-TextValue "D21" "4.098"
Set-TextValue "E21" "3.28%"
Set-TextValue "D22" "0.2088"
Set-TextValue "E22" "-0.80%"
Set-TextValue "D23" "0.04528"
Set-TextValue "E23" "-0.35%"
Set-TextValue "D24" "0.001217"
Set-TextValue "E24" "0.37%"
Set-TextValue "D25" "0.004359"
Set-TextValue "E25" "-6.12%"
Set-TextValue "D26" "0.0001293"
Set-TextValue "E26" "-0.62%"
Set-TextValue "D27" "0.0003375"
Set-TextValue "E27" "-0.78%"
Set-TextValue "D39" "0.01695"
Set-TextValue "E39" "5.92%"
Set-TextValue "D40" "0.04637"
Set-TextValue "E40" "4.38%"
Set-TextValue "D41" "0.007579"
Set-TextValue "E41" "3.25%"
Set-TextValue "D42" "0.009776"
Set-TextValue "E42" "13.25%"
Set-TextValue "D43" "0.1389"
Set-TextValue "E43" "4.31%"
Set-TextValue "D44" "0.002074"
Set-TextValue "E44" "2.51%"
Set-TextValue "D45" "0.009737"
Set-TextValue "E45" "2.19%"
Set-TextValue "D46" "0.00006065"
Set-TextValue "E46" "1.71%"
Set-TextValue "D47" "0.00000000746"
Set-TextValue "E47" "-0.79%"
Set-TextValue "D48" "2.551"
Set-TextValue "E48" "13.83%"
Set-TextValue "D49" "0.001989"
Set-TextValue "E49" "-31.40%"
Set-TextValue "D50" "0.00002088"
Set-TextValue "E50" "-0.79%"
Set-TextValue "D51" "0.0001988"
Set-TextValue "E51" "-0.79%"
